# Add RC filter components (R2, R23, R31 and C53, C57, C61, C67, C68, C69) to the
# telemetry string, regenerate outputs.
#
# This grows the "0.1uF / C6,C8,... " (SMD-0805C) designator group and shrinks
# the "0.1uF (SMD-1206C footprint)" group that previously absorbed
# C53/C57/C61, and inserts a relocated "1k" resistor row right after the
# "Value" row, shifting the other resistor rows down by one.
#
# Text-looking values are written with a leading apostrophe so Excel keeps
# them as text (matching the workbook's existing quoted-text cell style)
# instead of silently coercing look-alike numbers (e.g. "0.0", "60.4").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: 0.1uF / SMD-0805C group gains C53, C57, C61, C67, C68, C69 ---
$ws.Range("C5").Value = "C6, C8, C10, C11, C12, C17, C18, C27, C35, C37, C39, C45, C47, C49, C50, C51, C52, C53, C54, C55, C56, C57, C58, C59, C60, C61, C65, C66, C67, C68, C69"
$ws.Range("F5").Value = 31

# --- Row 7: 0.1uF / SMD-1206C group loses C53, C57, C61 (moved to row 5) ---
$ws.Range("C7").Value = "C20, C26"
$ws.Range("F7").Value = 2

# --- Rows 22-29: resistor table. A new "1k" row (R2, R23, R31, R40, R43, R46)
#     is inserted right after the "Value" row (row 21), pushing 100k/3.3k/
#     0.0/10k/60.4/249k/30k down by one row each. ---

$ws.Range("A22").Value = "'1k"
$ws.Range("B22").Value = "RES SMD 1K OHM 5% 1/8W 0805"
$ws.Range("C22").Value = "R2, R23, R31, R40, R43, R46"
$ws.Range("D22").Value = "SMD-0805-RES"
$ws.Range("E22").Value = "ERJ-6GEYJ102V"
$ws.Range("F22").Value = 6

$ws.Range("A23").Value = "'100k"
$ws.Range("B23").Value = "RES SMD 100K OHM 1% 1/8W 0805"
$ws.Range("C23").Value = "R3, R9, R24, R32"
$ws.Range("D23").Value = "SMD-0805-RES"
$ws.Range("E23").Value = "ERJ-6ENF1003V"
$ws.Range("F23").Value = 4

$ws.Range("A24").Value = "'3.3k"
$ws.Range("B24").Value = "RES SMD 3.3K OHM 1% 1/8W 0805"
$ws.Range("C24").Value = "R4, R10, R25, R33"
$ws.Range("D24").Value = "SMD-0805-RES"
$ws.Range("E24").Value = "ERJ-6ENF3301V"
$ws.Range("F24").Value = 4

$ws.Range("A25").Value = "'0.0"
$ws.Range("B25").Value = "RES SMD 0.0 OHM JUMPER 1/8W 0805"
$ws.Range("C25").Value = "R5, R6, R7, R17, R18, R19, R26, R27, R28, R34, R35, R36"
$ws.Range("D25").Value = "SMD-0805-RES"
$ws.Range("E25").Value = "ERJ-6GEY0R00V"
$ws.Range("F25").Value = 12

$ws.Range("A26").Value = "'10k"
$ws.Range("B26").Value = "RES SMD 10K OHM 5% 1/8W 0805, RES SMD 10K OHM 5% 1/2W 0805"
$ws.Range("C26").Value = "R8, R11, R12, R13, R16, R21, R29, R37, R38, R41, R44"
$ws.Range("D26").Value = "SMD-0805-RES"
$ws.Range("E26").Value = "ERJ-6GEYJ103V, ERJ-P06J103V"
$ws.Range("F26").Value = 11

$ws.Range("A27").Value = "'60.4"
$ws.Range("B27").Value = "RES SMD 60.4 OHM 1% 1/8W 0805"
$ws.Range("C27").Value = "R14, R15"
$ws.Range("D27").Value = "SMD-0805-RES"
$ws.Range("E27").Value = "ERJ-6ENF60R4V"
$ws.Range("F27").Value = 2

$ws.Range("A28").Value = "'249k"
$ws.Range("B28").Value = "RES SMD 249K OHM 1% 1/8W 0805"
$ws.Range("C28").Value = "R20"
$ws.Range("D28").Value = "SMD-0805-RES"
$ws.Range("E28").Value = "ERJ-6ENF2493V"
$ws.Range("F28").Value = 1

$ws.Range("A29").Value = "'30k"
$ws.Range("B29").Value = "RES SMD 30K OHM 1% 1/8W 0805"
$ws.Range("C29").Value = "R39, R42, R45"
$ws.Range("D29").Value = "SMD-0805-RES"
$ws.Range("E29").Value = "ERJ-6ENF3002V"
$ws.Range("F29").Value = 3
